# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    3  = 583
    5  = 297
    6  = 1114
    7  = 1450
    9  = 117
    12 = 177
    14 = 454
    15 = 1396
    16 = 127
    17 = 121
    22 = 1015
    26 = 6018
    31 = 14717
    32 = 1462
    33 = 234
    36 = 9505
    37 = 655
    39 = 164
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    3  = 583
    5  = 297
    6  = 1114
    7  = 1450
    9  = 117
    12 = 177
    14 = 454
    15 = 1396
    16 = 127
    17 = 121
    24 = 1015
    25 = 40
    29 = 6018
    34 = 14717
    35 = 1462
    36 = 234
    39 = 9505
    40 = 655
    42 = 164
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
